$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.813119411468506
$ws.Range("B1").Value = 4.69143009185791
$ws.Range("C1").Value = 3.841070890426636
$ws.Range("D1").Value = 0.9031198024749756
$ws.Range("E1").Value = 0.4741781949996948
